$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparação")

# Insert two new columns at F and G (old F,G,H shift to H,I,J)
$ws.Columns("F:G").Insert()

# Fix the "Concelho" label in row 2 (Setúbal -> Setubal)
$ws.Range("A2").Value = "Setubal"

# New header labels for the inserted columns
$ws.Range("F1").Value = "Ganho médio mensal Fem (2021)"
$ws.Range("G1").Value = "Ganho médio mensal Mas (2021)"

# Copy the currency number formatting (fill/border/numFmt) from column J
# onto the new E:G range for each data row, matching what Excel does when
# duplicating adjacent money-formatted cells.
$ws.Range("J2").Copy() | Out-Null
$ws.Range("E2:G2").PasteSpecial(-4122) | Out-Null

$ws.Range("J3").Copy() | Out-Null
$ws.Range("E3:G3").PasteSpecial(-4122) | Out-Null

$ws.Range("J4").Copy() | Out-Null
$ws.Range("E4:G4").PasteSpecial(-4122) | Out-Null

$ws.Range("J5").Copy() | Out-Null
$ws.Range("E5:G5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# New data values (average monthly earnings, female / male) per concelho
$ws.Range("F2").Value = 1087.8
$ws.Range("G2").Value = 1456.4

$ws.Range("F3").Value = 1242.7
$ws.Range("G3").Value = 1574.8

$ws.Range("F4").Value = 1019.6
$ws.Range("G4").Value = 1129.2

$ws.Range("F5").Value = 1006.8
$ws.Range("G5").Value = 1191.7

# Update current selection to match the saved workbook view
$ws.Range("G9").Select()
